$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference cell carrying the date number format already used by column D
$dateFormat = $ws.Cells.Item(110, 4).NumberFormat

$newRows = @(
    @{ Row = 111; A = 3; B = "Femacal de La Calera"; C = "Coquimbo"; D = 44509; E = 5; F = "Fruta"; G = 100107; H = "Otros"; I = 100107002; J = "Chirimoya"; K = "Cultivar IV Región"; L = "Especial"; M = 45; N = 27000; O = 27000; P = 27000; Q = "`$/bandeja 10 kilos"; R = "Provincia de Limarí"; S = 2700; T = 10 },
    @{ Row = 112; A = 3; B = "Femacal de La Calera"; C = "Coquimbo"; D = 44509; E = 5; F = "Fruta"; G = 100107; H = "Otros"; I = 100107002; J = "Chirimoya"; K = "Cultivar IV Región"; L = "Primera"; M = 48; N = 25000; O = 25000; P = 25000; Q = "`$/bandeja 10 kilos"; R = "Provincia de Limarí"; S = 2500; T = 10 },
    @{ Row = 113; A = 3; B = "Femacal de La Calera"; C = "Coquimbo"; D = 44509; E = 5; F = "Fruta"; G = 100107; H = "Otros"; I = 100107002; J = "Chirimoya"; K = "Cultivar IV Región"; L = "Segunda"; M = 40; N = 22000; O = 22000; P = 22000; Q = "`$/bandeja 10 kilos"; R = "Provincia de Limarí"; S = 2200; T = 10 }
)

foreach ($rowData in $newRows) {
    $r = $rowData.Row

    $ws.Cells.Item($r, 1).Value2 = $rowData.A
    $ws.Cells.Item($r, 2).Value2 = $rowData.B
    $ws.Cells.Item($r, 3).Value2 = $rowData.C

    $dCell = $ws.Cells.Item($r, 4)
    $dCell.Value2 = $rowData.D
    $dCell.NumberFormat = $dateFormat

    $ws.Cells.Item($r, 5).Value2 = $rowData.E
    $ws.Cells.Item($r, 6).Value2 = $rowData.F
    $ws.Cells.Item($r, 7).Value2 = $rowData.G
    $ws.Cells.Item($r, 8).Value2 = $rowData.H
    $ws.Cells.Item($r, 9).Value2 = $rowData.I
    $ws.Cells.Item($r, 10).Value2 = $rowData.J
    $ws.Cells.Item($r, 11).Value2 = $rowData.K
    $ws.Cells.Item($r, 12).Value2 = $rowData.L
    $ws.Cells.Item($r, 13).Value2 = $rowData.M
    $ws.Cells.Item($r, 14).Value2 = $rowData.N
    $ws.Cells.Item($r, 15).Value2 = $rowData.O
    $ws.Cells.Item($r, 16).Value2 = $rowData.P
    $ws.Cells.Item($r, 17).Value2 = $rowData.Q
    $ws.Cells.Item($r, 18).Value2 = $rowData.R
    $ws.Cells.Item($r, 19).Value2 = $rowData.S
    $ws.Cells.Item($r, 20).Value2 = $rowData.T
}
